$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 347
$ws.Range("F6").Value = 521
$ws.Range("F8").Value = 237
$ws.Range("F12").Value = 72

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 3

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F10").Value = 9233
$ws.Range("F11").Value = 237
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 641
$ws.Range("F17").Value = 0
